$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 708; rows 708..749 shift down to 709..750
$ws.Rows.Item(708).Insert()

# The date in column A is stored as text (matching the rest of the sheet),
# so copy it (values-only) from the row above rather than re-typing it,
# which keeps Excel from "helpfully" re-interpreting it as a real date.
$ws.Range("A707").Copy()
$ws.Range("A708").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(708, 2).Value = "火"
$ws.Cells.Item(708, 3).Value = 19
$ws.Cells.Item(708, 4).Value = 30
